$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '45.083.01'
$ws.Range("E2").Value = '  +3.58%  '
$ws.Range("D3").Value = '2.426.63'
$ws.Range("E3").Value = '  +0.72%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '317.42'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +3.55%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '102.36'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +5.47%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.516'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +1.38%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("E9").Value = '  +7.04%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.33'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +0.76%  '
$ws.Range("E11").Value = '  +0.74%  '
$ws.Range("E12").Value = '  -2.49%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.12'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -1.97%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.00'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +1.56%  '
$ws.Range("D15").Value = '2.806.79'
$ws.Range("E15").Value = '  +0.98%  '
$ws.Range("D16").Value = '2.423.74'
$ws.Range("E16").Value = '  +0.53%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.839'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +2.01%  '
$ws.Range("D18").Value = '45.010.35'
$ws.Range("E18").Value = '  +3.48%  '
$ws.Range("E19").Value = '  +0.80%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.35'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -1.10%  '
$ws.Range("E21").Value = '  +2.27%  '
$ws.Range("E22").Value = '  +0.56%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '243.92'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +2.49%  '
$ws.Range("E24").Value = '  +0.57%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.49'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.63%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.32'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +1.30%  '
$ws.Range("B28").Value = 'Cosmos'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.55'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +1.36%  '
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.07'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -11.94%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '49.22'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +2.80%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '32.82'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +1.58%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.15'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +9.46%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.124'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +7.42%  '
$ws.Range("E34").Value = '  +1.69%  '
$ws.Range("E35").Value = '  +0.23%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0765'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +2.05%  '
$ws.Range("E37").Value = '  -0.43%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.42'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +0.73%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.86'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -1.91%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '124.52'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -7.14%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.21'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -2.63%  '
$ws.Range("E42").Value = '  +0.62%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '20.76'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -2.82%  '
$ws.Range("E44").Value = '  +1.83%  '
$ws.Range("D45").Value = '1.932.48'
$ws.Range("E45").Value = '  -0.74%  '
$ws.Range("E46").Value = '  -3.12%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.93'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +3.95%  '
$ws.Range("B48").Value = 'FraxShare'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.25'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -0.10%  '
$ws.Range("B49").Value = 'Stacks'
$ws.Range("C49").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.80'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +15.79%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '76.57'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +5.95%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '53.83'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +2.14%  '
